# Update the Epcam-Epcam NATMI LR-pairs sheet with the new TPM-derived
# values: the single ECs/FAPs self-pair row is replaced by the full
# 2x2 cross of Sending/Target cluster (ECs, FAPs) combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target
# cluster, E..T the NATMI expression / specificity metrics.
$rows = @(
    @("ECs","Epcam","Epcam","ECs",
      1,0.5,0.0582165,0.116433,0.4693547931970906,0.3709368598453593,
      1,0.5,0.0582165,0.116433,0.4693547931970906,0.3709368598453593,
      0.00338916087225,0.013556643489,0.2202939218970836,0.1375941539919358),
    @("ECs","Epcam","Epcam","FAPs",
      1,0.5,0.0582165,0.116433,0.4693547931970906,0.3709368598453593,
      1,0.3333333333333333,0.06581866666666666,0.197456,0.5306452068029094,0.6290631401546407,
      0.003831732408,0.022990394448,0.2490608713000069,0.2333427058534236),
    @("FAPs","Epcam","Epcam","ECs",
      1,0.3333333333333333,0.06581866666666666,0.197456,0.5306452068029094,0.6290631401546407,
      1,0.5,0.0582165,0.116433,0.4693547931970906,0.3709368598453593,
      0.003831732408,0.022990394448,0.2490608713000069,0.2333427058534236),
    @("FAPs","Epcam","Epcam","FAPs",
      1,0.3333333333333333,0.06581866666666666,0.197456,0.5306452068029094,0.6290631401546407,
      1,0.3333333333333333,0.06581866666666666,0.197456,0.5306452068029094,0.6290631401546407,
      0.004332096881777777,0.03898887193599999,0.2815843355029025,0.3957204343012171)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowData[$c]
    }
}
